$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to write. Values that Excel would otherwise auto-convert
# to numbers (plain decimal-looking strings) are prefixed with a leading
# apostrophe so Excel stores them as literal text, matching the source data
# (coinranking price/volume strings are text, not numbers).
$updates = [ordered]@{
    "D2" = "65.533.47"
    "E2" = "  -0.69%  "
    "D3" = "3.337.29"
    "E3" = "  -4.07%  "
    "E4" = "  +0.02%  "
    "D5" = "'574.88"
    "E5" = "  -1.29%  "
    "D6" = "'178.89"
    "E6" = "  +3.12%  "
    "D7" = "'0.620"
    "E7" = "  +3.63%  "
    "E8" = "  +0.01%  "
    "D9" = "3.335.30"
    "E9" = "  -4.09%  "
    "E10" = "  -1.76%  "
    "D11" = "'6.86"
    "E11" = "  -0.49%  "
    "D12" = "'0.406"
    "E12" = "  -0.53%  "
    "D13" = "3.919.70"
    "E13" = "  -3.99%  "
    "D14" = "'0.135"
    "E14" = "  +0.37%  "
    "D15" = "'28.39"
    "E15" = "  -4.28%  "
    "D16" = "65.523.22"
    "E16" = "  -0.78%  "
    "E17" = "  -1.60%  "
    "D18" = "3.335.63"
    "E18" = "  -4.17%  "
    "D19" = "'5.76"
    "E19" = "  -2.72%  "
    "D20" = "'13.42"
    "E20" = "  -3.31%  "
    "D21" = "'364.67"
    "E21" = "  -0.58%  "
    "E22" = "  -3.91%  "
    "E23" = "  +0.13%  "
    "D24" = "'71.49"
    "E24" = "  -1.91%  "
    "E25" = "  -3.10%  "
    "D26" = "'0.518"
    "E26" = "  -3.05%  "
    "D27" = "'9.53"
    "E27" = "  -1.50%  "
    "E28" = "  -0.79%  "
    "E29" = "  -0.02%  "
    "E30" = "  -1.52%  "
    "D31" = "'5.63"
    "E31" = "  -1.93%  "
    "E32" = "  +0.00%  "
    "D33" = "'22.92"
    "E33" = "  -4.15%  "
    "D34" = "'6.80"
    "E34" = "  -4.65%  "
    "E35" = "  -6.08%  "
    "E36" = "  -2.75%  "
    "D37" = "'159.91"
    "E37" = "  -0.41%  "
    "E38" = "  -4.92%  "
    "D39" = "'27.32"
    "E39" = "  -6.80%  "
    "E40" = "  -0.51%  "
    "D41" = "2.721.19"
    "E41" = "  -2.88%  "
    "D42" = "'2.52"
    "E42" = "  -3.81%  "
    "D43" = "'6.23"
    "E43" = "  -3.72%  "
    "D44" = "'4.27"
    "E44" = "  -4.24%  "
    "D45" = "'39.93"
    "E45" = "  -0.49%  "
    "D46" = "'0.0666"
    "E46" = "  -2.49%  "
    "D47" = "'331.33"
    "E47" = "  +1.53%  "
    "D48" = "'23.99"
    "E48" = "  -0.64%  "
    "E49" = "  -3.92%  "
    "E50" = "  +2.62%  "
    "B51" = "Arweave"
    "C51" = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
    "D51" = "'30.39"
    "E51" = "  -1.07%  "
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
